# Updated cryptos list - price/volume refresh + minor ranking reshuffle (rows 36-37, 39-40, 47-48)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextCell "D2" "63.944.07"
Set-TextCell "E2" "  -4.05%  "

Set-TextCell "D3" "3.135.69"
Set-TextCell "E3" "  -2.98%  "

Set-TextCell "E4" "  +0.04%  "

Set-TextCell "D5" "601.30"
Set-TextCell "E5" "  -0.52%  "

Set-TextCell "D6" "145.64"
Set-TextCell "E6" "  -7.54%  "

Set-TextCell "E7" "  +0.13%  "

Set-TextCell "D8" "3.135.43"
Set-TextCell "E8" "  -2.91%  "

Set-TextCell "D9" "0.521"
Set-TextCell "E9" "  -4.10%  "

Set-TextCell "D10" "0.150"
Set-TextCell "E10" "  -6.96%  "

Set-TextCell "D11" "5.49"
Set-TextCell "E11" "  -4.55%  "

Set-TextCell "D12" "0.471"
Set-TextCell "E12" "  -6.27%  "

Set-TextCell "D13" "0.0000253"
Set-TextCell "E13" "  -6.32%  "

Set-TextCell "D14" "35.81"
Set-TextCell "E14" "  -8.29%  "

Set-TextCell "D15" "3.663.46"
Set-TextCell "E15" "  -2.68%  "

Set-TextCell "D16" "64.012.32"
Set-TextCell "E16" "  -4.03%  "

Set-TextCell "E17" "  +0.60%  "

Set-TextCell "D18" "3.145.20"
Set-TextCell "E18" "  -2.85%  "

Set-TextCell "D19" "6.88"
Set-TextCell "E19" "  -5.66%  "

Set-TextCell "D20" "477.11"
Set-TextCell "E20" "  -6.05%  "

Set-TextCell "D21" "14.43"
Set-TextCell "E21" "  -5.75%  "

Set-TextCell "D22" "0.701"
Set-TextCell "E22" "  -5.10%  "

Set-TextCell "D23" "7.67"
Set-TextCell "E23" "  -4.77%  "

Set-TextCell "D24" "13.66"
Set-TextCell "E24" "  -6.64%  "

Set-TextCell "D25" "82.36"
Set-TextCell "E25" "  -4.06%  "

Set-TextCell "E26" "  +0.00%  "

Set-TextCell "D27" "2.87"
Set-TextCell "E27" "  -4.72%  "

Set-TextCell "D28" "8.38"
Set-TextCell "E28" "  -7.60%  "

Set-TextCell "D29" "2.19"
Set-TextCell "E29" "  -7.11%  "

Set-TextCell "D30" "0.118"
Set-TextCell "E30" "  -27.56%  "

Set-TextCell "D31" "6.79"
Set-TextCell "E31" "  -3.02%  "

Set-TextCell "D32" "2.73"
Set-TextCell "E32" "  -6.44%  "

Set-TextCell "E33" "  +0.06%  "

Set-TextCell "D34" "26.05"
Set-TextCell "E34" "  -7.82%  "

Set-TextCell "D35" "1.10"
Set-TextCell "E35" "  -5.69%  "

Set-TextCell "B36" "OKB"
Set-TextCell "C36" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D36" "54.16"
Set-TextCell "E36" "  -2.13%  "

Set-TextCell "B37" "Filecoin"
Set-TextCell "C37" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D37" "5.97"
Set-TextCell "E37" "  -6.99%  "

Set-TextCell "D38" "0.0₃0724"
Set-TextCell "E38" "  -7.87%  "

Set-TextCell "B39" "Bittensor"
Set-TextCell "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D39" "450.26"
Set-TextCell "E39" "  -9.83%  "

Set-TextCell "B40" "dogwifhat"
Set-TextCell "C40" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D40" "2.95"
Set-TextCell "E40" "  -8.08%  "

Set-TextCell "D41" "0.0394"
Set-TextCell "E41" "  -6.42%  "

Set-TextCell "D42" "0.119"
Set-TextCell "E42" "  -6.98%  "

Set-TextCell "D43" "8.38"
Set-TextCell "E43" "  -4.01%  "

Set-TextCell "D44" "2.846.83"
Set-TextCell "E44" "  -3.59%  "

Set-TextCell "D45" "0.267"
Set-TextCell "E45" "  -9.41%  "

Set-TextCell "D46" "2.25"
Set-TextCell "E46" "  -8.50%  "

Set-TextCell "B47" "USDe"
Set-TextCell "C47" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D47" "0.998"
Set-TextCell "E47" "  -0.07%  "

Set-TextCell "B48" "InjectiveProtocol"
Set-TextCell "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D48" "26.08"
Set-TextCell "E48" "  -7.96%  "

Set-TextCell "D49" "2.30"
Set-TextCell "E49" "  -5.33%  "

Set-TextCell "D50" "0.114"
Set-TextCell "E50" "  -3.81%  "

Set-TextCell "D51" "117.17"
Set-TextCell "E51" "  -3.66%  "

